$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 176445
$ws.Range("C4").Value = 166402
$ws.Range("C7").Value = 5.69
$ws.Range("C8").Value = 64.65000000000001
